$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 40
$ws.Range("A40").Value = 42036
$ws.Range("A40").NumberFormat = $ws.Range("A39").NumberFormat
$ws.Range("B40").Value = 51.73
$ws.Range("C40").Value = 17
$ws.Range("D40").Value = 34.729999999999997
$ws.Range("H40").Value = 1.1200000000000001
$ws.Range("I40").Value = 18.98
$ws.Range("J40").Value = 7.89
$ws.Range("K40").Value = 6.65
$ws.Range("M40").Value = 68

# Row 41
$ws.Range("A41").Value = 42064
$ws.Range("A41").NumberFormat = $ws.Range("A39").NumberFormat
$ws.Range("B41").Value = 53.32
$ws.Range("C41").Value = 17
$ws.Range("D41").Value = 36.32
$ws.Range("H41").Value = 3.6
$ws.Range("I41").Value = 32.72
$ws.Range("M41").Value = 71

# Update selection to match diff (activeCell Q32)
$ws.Range("Q32").Select()
